$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Median Value" (C) and "Tier" (D) columns with scores recalculated
# relative to the median AFTER merging with zip/census tract data.
$data = @{
  2 = @(0.6367061895422116, "Below Median")
  3 = @(2.151694113166132, "1st Tier")
  4 = @(0.4152559601949424, "Below Median")
  5 = @(1.480351279544436, "2nd Tier")
  6 = @(1.816164731501566, "1st Tier")
  7 = @(0.6751159698831268, "Below Median")
  8 = @(0.8343256418679427, "Below Median")
  9 = @(0.8495813633245719, "Below Median")
  10 = @(1.662427194835814, "2nd Tier")
  11 = @(1.574504923760461, "2nd Tier")
  12 = @(0.7008685365784988, "Below Median")
  13 = @(1.051779167353432, "4th Tier")
  14 = @(1.578132928486704, "2nd Tier")
  15 = @(0.9660305703137798, "Below Median")
  16 = @(0.9262786326610984, "Below Median")
  17 = @(0.7113661375202767, "Below Median")
  18 = @(0.5943664124520932, "Below Median")
  19 = @(0.9706532156137764, "Below Median")
  20 = @(1.971854147155835, "1st Tier")
  21 = @(1.029056510453705, "4th Tier")
  22 = @(0.6041801011699841, "Below Median")
  23 = @(0.7042547396290442, "Below Median")
  24 = @(1.58090015259156, "2nd Tier")
  25 = @(1.274431542096313, "3rd Tier")
  26 = @(1.017195898095257, "4th Tier")
  27 = @(0.9527646380022146, "Below Median")
  28 = @(0.7544713918411923, "Below Median")
  29 = @(1.810803237471412, "1st Tier")
  30 = @(1.114186473952443, "4th Tier")
  31 = @(0.5220518108279673, "Below Median")
  32 = @(1.271786279135188, "3rd Tier")
  33 = @(0.803645219242051, "Below Median")
  34 = @(0.8064669670138866, "Below Median")
  35 = @(1.285297181325016, "3rd Tier")
  36 = @(0.972578313468646, "Below Median")
  37 = @(1.092784966663571, "4th Tier")
  38 = @(1.332853528599145, "3rd Tier")
  39 = @(0.5161523701148558, "Below Median")
  40 = @(1.249578634798395, "3rd Tier")
  41 = @(0.5775327841399486, "Below Median")
  42 = @(0.8728116804391033, "Below Median")
  43 = @(0.9603021918282103, "Below Median")
  44 = @(1.792195576980931, "1st Tier")
  45 = @(1.328054682765954, "3rd Tier")
  46 = @(1.225163696714938, "3rd Tier")
  47 = @(2.508482521899435, "1st Tier")
  48 = @(0.8090840305675493, "Below Median")
  49 = @(1.185328245786839, "3rd Tier")
  50 = @(1.689946331205805, "2nd Tier")
  51 = @(1.061578259065911, "4th Tier")
  52 = @(0.4005772878787622, "Below Median")
  53 = @(1.752304306788031, "1st Tier")
  54 = @(0.3908166912043015, "Below Median")
  55 = @(1.281799565195446, "3rd Tier")
  56 = @(1.055425278245005, "4th Tier")
  57 = @(0.7133060869313723, "Below Median")
  58 = @(0.9082139757205137, "Below Median")
  59 = @(1.164488941466887, "4th Tier")
  60 = @(0.975475692618172, "Below Median")
  61 = @(1.093860511706489, "4th Tier")
  62 = @(1.003582302218997, "4th Tier")
  63 = @(0.8837556744572219, "Below Median")
  64 = @(1.847294697855373, "1st Tier")
  65 = @(0.8605934767954495, "Below Median")
  66 = @(0.6550104486463219, "Below Median")
  67 = @(1.376570428220897, "2nd Tier")
  68 = @(0.7362833121018495, "Below Median")
  69 = @(0.8218236205635329, "Below Median")
  70 = @(0.811746439279329, "Below Median")
  71 = @(0.8407651782192986, "Below Median")
  72 = @(1.694656653543663, "2nd Tier")
  73 = @(0.7232782312579946, "Below Median")
  74 = @(1.787881131648522, "1st Tier")
  75 = @(2.337542246996575, "1st Tier")
  76 = @(1.538267214068107, "2nd Tier")
  77 = @(0.3657798483147726, "Below Median")
  78 = @(0.837845077440315, "Below Median")
  79 = @(1.320452205134144, "3rd Tier")
  80 = @(1.615973575159628, "2nd Tier")
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 3).Value = $vals[0]
  $ws.Cells.Item($row, 4).Value = $vals[1]
}
